# Updated cryptos list on Fri Nov 17 14:54:56 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures on the
# cryptos worksheet. Rows 46 and 47 also swap rank order (ARBITRUM moves
# up to rank 44 / row 46, FraxShare moves down to rank 45 / row 47), so
# those two rows get their Coin name, Link and Price/Volume fully
# replaced instead of just the numbers.
#
# Note: several Price values (e.g. "241.19", "0.602") look like plain
# numbers, but in the source file column D is text. Assigning such a
# string straight to .Value would make Excel's COM layer silently
# reinterpret it as a number, so for those cells we briefly force the
# cell to Text format, assign the literal string, then drop the style
# back to Normal so no stray number format is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.101.24'
$ws.Range('E2').Value = '  -1.93%  '
$ws.Range('D3').Value = '1.932.67'
$ws.Range('E3').Value = '  -4.78%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.19'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.602'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.77%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '56.55'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -10.39%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.363'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '55.07'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0815'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.03%  '
$ws.Range('E12').Value = '  -0.53%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.818'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -8.96%  '
$ws.Range('D14').Value = '2.226.59'
$ws.Range('E14').Value = '  -4.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.13'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -9.96%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.21'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -8.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.19'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -6.66%  '
$ws.Range('D18').Value = '1.936.41'
$ws.Range('E18').Value = '  -4.38%  '
$ws.Range('D19').Value = '36.095.54'
$ws.Range('E19').Value = '  -1.87%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.35'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.19%  '
$ws.Range('D21').Value = '0.0₃0857'
$ws.Range('E21').Value = '  -3.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '227.15'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.94'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -8.33%  '
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.47'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.80%  '
$ws.Range('E26').Value = '  -2.34%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.30'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.68'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.17'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -6.15%  '
$ws.Range('E30').Value = '  -15.02%  '
$ws.Range('E31').Value = '  -3.55%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.13'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.67%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.65'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -8.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0622'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.25'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.48%  '
$ws.Range('E36').Value = '  +0.36%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.83%  '
$ws.Range('E38').Value = '  -2.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.12'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -11.47%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.81'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -14.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0962'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.47%  '
$ws.Range('E42').Value = '  -2.35%  '
$ws.Range('E43').Value = '  -7.38%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0207'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '15.41'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -9.86%  '

# Row 46 / 47: ARBITRUM and FraxShare swap rank order.
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.02'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -10.02%  '

$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.27'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.35%  '

$ws.Range('D48').Value = '1.326.27'
$ws.Range('E48').Value = '  -3.37%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '86.90'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.81'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '45.77'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.53%  '
